$d = $word.ActiveDocument

# Paragraph 1 (Title style): merge the word-by-word runs of
# "Questions: Rationalizing the denominator" into a single run.
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
[void]$r1.MoveEnd(1, -1)
[void]$r1.Find.Execute("Questions: Rationalizing the denominator", $false, $false, $false, $false, $false, $true, 1, $false, "Questions: Rationalizing the denominator", 2)

# Paragraph 2 (Author style): merge "Maximilian" / " " / "Volmar" into one run.
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
[void]$r2.MoveEnd(1, -1)
[void]$r2.Find.Execute("Maximilian Volmar", $false, $false, $false, $false, $false, $true, 1, $false, "Maximilian Volmar", 2)

# Paragraph 4 (Abstract style): merge the long word-by-word sentence into one run.
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
[void]$r4.MoveEnd(1, -1)
[void]$r4.Find.Execute("A selection of questions for the study guide on rationalizing the denominator.", $false, $false, $false, $false, $false, $true, 1, $false, "A selection of questions for the study guide on rationalizing the denominator.", 2)
